# Update timestamps (column A) and actual production values (column B)
# for rows 2..97 per the commit: shift dates forward by one day and
# refresh production readings (issue on PC Sun model requiring retraining).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(45970.01041666666,45970.02083333334,45970.03125,45970.04166666666,45970.05208333334,45970.0625,45970.07291666666,45970.08333333334,45970.09375,45970.10416666666,45970.11458333334,45970.125,45970.13541666666,45970.14583333334,45970.15625,45970.16666666666,45970.17708333334,45970.1875,45970.19791666666,45970.20833333334,45970.21875,45970.22916666666,45970.23958333334,45970.25,45970.26041666666,45970.27083333334,45970.28125,45970.29166666666,45970.30208333334,45970.3125,45970.32291666666,45970.33333333334,45970.34375,45970.35416666666,45970.36458333334,45970.375,45970.38541666666,45970.39583333334,45970.40625,45970.41666666666,45970.42708333334,45970.4375,45970.44791666666,45970.45833333334,45970.46875,45970.47916666666,45970.48958333334,45970.5,45970.51041666666,45970.52083333334,45970.53125,45970.54166666666,45970.55208333334,45970.5625,45970.57291666666,45970.58333333334,45970.59375,45970.60416666666,45970.61458333334,45970.625,45970.63541666666,45970.64583333334,45970.65625,45970.66666666666,45970.67708333334,45970.6875,45970.69791666666,45970.70833333334,45970.71875,45970.72916666666,45970.73958333334,45970.75,45970.76041666666,45970.77083333334,45970.78125,45970.79166666666,45970.80208333334,45970.8125,45970.82291666666,45970.83333333334,45970.84375,45970.85416666666,45970.86458333334,45970.875,45970.88541666666,45970.89583333334,45970.90625,45970.91666666666,45970.92708333334,45970.9375,45970.94791666666,45970.95833333334,45970.96875,45970.97916666666,45970.98958333334,45971)
$bValues = @(1400,1404,1345,1328,1350,1336,1301,1319,1391,1318,1257,1263,1210,1188,1201,1193,1150,1111,1077,1062,1045,1012,992,970,928,903,895,847,783,703,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
